# Generate Report for Handoff
# Replace the first file's (3bc449a5-...) and second file's (a74d51bb-...)
# localization-status rows with the new handoff data for
# e3269309-456d-45f8-a05d-1248095a1f4b.md and ffff01fc7a75-bd34-4f19-90b6-2cb650ec98d2.md.

$wb = $excel.ActiveWorkbook

$newFile1 = "e3269309-456d-45f8-a05d-1248095a1f4b.md"
$newFile2 = "ffff01fc7a75-bd34-4f19-90b6-2cb650ec98d2.md"
$status   = "Ready for handoff"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value() = $newFile1
$wsZh.Range("C2").Value() = $status
$wsZh.Range("G2").Value() = "e3269309-456d-45f8-a05d-1248095a1f4b.0fc09a7c65051cc44f5653842e410e625810c261.zh-cn.xlf"
$wsZh.Range("H2").Value() = "2016-08-28 00:59:47"
$wsZh.Range("I2").Value() = ""
$wsZh.Range("J2").Value() = ""
$wsZh.Range("K2").Value() = "0001-01-01 00:00:00"

$wsZh.Range("A3").Value() = $newFile2
$wsZh.Range("C3").Value() = $status
$wsZh.Range("F3").Value() = "True"
$wsZh.Range("G3").Value() = "e3269309-456d-45f8-a05d-1248095a1f4b.0fc09a7c65051cc44f5653842e410e625810c261.zh-cn.xlf"
$wsZh.Range("H3").Value() = "2016-08-28 00:59:47"
$wsZh.Range("I3").Value() = ""
$wsZh.Range("J3").Value() = ""
$wsZh.Range("K3").Value() = "0001-01-01 00:00:00"

$wsZh.Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsZh.Hyperlinks.Item(2).TextToDisplay = $newFile2

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value() = $newFile1
$wsDe.Range("C2").Value() = $status
$wsDe.Range("G2").Value() = "e3269309-456d-45f8-a05d-1248095a1f4b.0fc09a7c65051cc44f5653842e410e625810c261.de-de.xlf"
$wsDe.Range("H2").Value() = "2016-08-28 00:59:52"
$wsDe.Range("I2").Value() = ""
$wsDe.Range("J2").Value() = ""
$wsDe.Range("K2").Value() = "0001-01-01 00:00:00"

$wsDe.Range("A3").Value() = $newFile2
$wsDe.Range("C3").Value() = $status
$wsDe.Range("F3").Value() = "True"
$wsDe.Range("G3").Value() = "e3269309-456d-45f8-a05d-1248095a1f4b.0fc09a7c65051cc44f5653842e410e625810c261.de-de.xlf"
$wsDe.Range("H3").Value() = "2016-08-28 00:59:52"
$wsDe.Range("I3").Value() = ""
$wsDe.Range("J3").Value() = ""
$wsDe.Range("K3").Value() = "0001-01-01 00:00:00"

$wsDe.Hyperlinks.Item(1).TextToDisplay = $newFile1
$wsDe.Hyperlinks.Item(2).TextToDisplay = $newFile2

# ---- Overview sheet ----
# Cell values here reference the same shared strings as the detail sheets
# (File Name / Path And Name / zh-cn / de-de / Latest HO date), so they update
# automatically. Only refresh the hyperlink display text.
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Hyperlinks.Item(1).TextToDisplay = "e2e\" + $newFile1
$wsOv.Hyperlinks.Item(2).TextToDisplay = "e2e\" + $newFile2
